# Adds team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 45

# New header cells in row 1: reuse the formatting of the existing header
# cells (bold, bordered, centered) by copy/paste-special of an existing
# header cell's format, then set the text values.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2..45: same W/L/T record repeated for every player row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 78   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 84   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
